# 2023 - Updated timing day 14 part II
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Day 14 (row 20): add Part I / Part II timings
$ws.Range("D20").Value = 3
$ws.Range("E20").Value = 7

# Recalculate formulas (G20, G32, D32, E32, G34, D34, E34, etc.)
$excel.Calculate()

# Update the view: active selection, matching the author's final
# on-screen state when the workbook was saved.
$ws.Range("E21").Select()
